$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.233.35"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.867.00"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7207"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07803"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3085"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08262"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "1.901.18"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7207"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.222"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").Value = "29.284.97"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.856"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007800"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.121.38"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.968"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1605"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.935"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.343"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.493"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.395"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.089"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05198"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.926"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.182"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7276"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.681"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01850"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.703"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "1.167.93"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9038"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.111"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.021.89"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5290"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.779"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("E49").Value = "  +4.27%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.309"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("B51").Value = "SynthetixNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.874"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.61%  "
